# Applies price/volume updates to the cryptos worksheet (rows 2-51, columns D & E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.347.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.550.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.550.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.153.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.538.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.434.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "430.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.690.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.544.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("E35").Value = "  -5.18%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0846"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.888"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  +4.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.23%  "
